# Apply the edits described by the commit:
# - remove the '__groups__' worksheet
# - update a couple of population figures on the 'pop' sheet
# - leave '__axes__' as the active/selected sheet

$wb = $excel.ActiveWorkbook

# Update changed population values on sheet "pop"
$wsPop = $wb.Worksheets.Item("pop")
$wsPop.Range("D4").Value = 32045129
$wsPop.Range("E4").Value = 32174258
$wsPop.Range("D5").Value = 34120851
$wsPop.Range("E5").Value = 34283895

# Remove the '__groups__' sheet entirely
$excel.DisplayAlerts = $false
$wsGroups = $wb.Worksheets.Item("__groups__")
$wsGroups.Delete()
$excel.DisplayAlerts = $true

# Make '__axes__' the active sheet (it becomes the last remaining sheet)
$wsAxes = $wb.Worksheets.Item("__axes__")
$wsAxes.Activate()
